$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.832.30"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "1.561.09"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "205.44"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "21.57"
$ws.Range("E8").Value = "  -3.08%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "1.781.44"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "1.571.33"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("D16").Value = "26.835.08"
$ws.Range("E16").Value = "  -0.99%  "
$ws.Range("D17").Value = "61.28"
$ws.Range("E17").Value = "  -2.61%  "
$ws.Range("D18").Value = "215.28"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "7.33"
$ws.Range("E19").Value = "  +1.50%  "
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "4.12"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("D25").Value = "153.36"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D26").Value = "6.61"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("E31").Value = "  -3.47%  "
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").Value = "1.381.27"
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("E37").Value = "  -2.23%  "
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("E39").Value = "  +1.85%  "
$ws.Range("D40").Value = "0.809"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("E43").Value = "  +4.82%  "
$ws.Range("E44").Value = "  -0.40%  "
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.50"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "1.695.67"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").Value = "86.49"
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("D49").Value = "0.0508"
$ws.Range("E49").Value = "  +3.14%  "
$ws.Range("D50").Value = "0.0₇0968"
$ws.Range("E50").Value = "  -2.54%  "
$ws.Range("E51").Value = "  +0.84%  "
